# Final 2023 generation edits for the 4CEMERLANG homeroom statement.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Theme Party" activity (row 40, column C) to the actual
#    event that took place.
$ws.Range("C40").Value = "Kolokium Zon Selatan"

# 2. Record the Demerit amount for the "Semakan Kali Ketiga" line under the
#    "Penyertaan Pertandingan" section (row 24).
$ws.Range("E24").Value = 400

# 3. "Kolokium Zon Selatan" (row 40) did not actually earn the previously
#    estimated 100 merit points - correct it to 0.
$ws.Range("D40").Value = 0

# 4. Merge B15:C15 (the "Penandaan Fail" label row), matching the merge
#    pattern already used by the surrounding rows (B13:C13, B14:C14).
$ws.Range("B15:C15").Merge()

# Merging leaves C15 without any content, so drop its leftover bold/wrap
# formatting back to the plain default font while keeping the cell's
# border, matching the final look of the other merged label cells.
$ws.Range("C15").Font.Name = "Calibri"
$ws.Range("C15").Font.Size = 11
$ws.Range("C15").WrapText = $false
$ws.Range("C15").ShrinkToFit = $false
